# Generate Report for Handback
# ------------------------------------------------------------------
# This script reproduces a "handback" localization-status report run:
#   * the Overview status text moves from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   * each language sheet (zh-cn / de-de) gets its "Latest Target File"
#     (col I) and "Latest Handback File" (col J) populated, with col I
#     turned into a hyperlink to the handed-back markdown file, and the
#     "Latest Handback DateTime" (col K) stamped with the handback time.
#   * a few columns are widened so the new, longer values are readable.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdbee8bb0cadca47eacfdb25318c0ba91d04e6b9/e2e"

$fileA = "ac55d992-1edd-410e-bdfd-f7bebed1963e.md"
$fileB = "c9c70c95-c479-43fd-a4d2-48270d45b7c0.md"

# ------------------------------------------------------------------
# Overview sheet: status text + column widths
# ------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

$ovw.Columns.Item(5).ColumnWidth = 29.17
$ovw.Columns.Item(6).ColumnWidth = 29.17

# ------------------------------------------------------------------
# zh-cn sheet
# ------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Columns.Item(3).ColumnWidth = 29.17
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

# "Status" column shares the same text as the Overview status cells.
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 (ac55d992-...)
$zh.Range("I2").Value = $fileA
$zh.Range("J2").Value = "ac55d992-1edd-410e-bdfd-f7bebed1963e.345e5284197f96f7f8137abb0f1f643e2aedff4e.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-24 07:03:05"

# Row 3 (c9c70c95-...)
$zh.Range("I3").Value = $fileB
$zh.Range("J3").Value = "c9c70c95-c479-43fd-a4d2-48270d45b7c0.5ce180b0445c7c4ad289ac1955311cb07e44f947.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-24 07:03:05"

# Rebuild hyperlinks so the new Latest-Target-File links land right after
# each row's source-file link (A2, I2, A3, I3), matching relationship order.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "$repoBase/$fileA", "", "", $fileA)
$zh.Hyperlinks.Add($zh.Range("I2"), "$repoBase/$fileA", "", "", $fileA)
$zh.Hyperlinks.Add($zh.Range("A3"), "$repoBase/$fileB", "", "", $fileB)
$zh.Hyperlinks.Add($zh.Range("I3"), "$repoBase/$fileB", "", "", $fileB)

# Match the underline/blue "HyperLink" look used by the A-column links.
$zh.Range("I2").Font.Underline = 1
$zh.Range("I2").Font.Name = "Calibri"
$zh.Range("I2").Font.Size = 11
$zh.Range("I2").Font.Color = 15570276
$zh.Range("I3").Font.Underline = 1
$zh.Range("I3").Font.Name = "Calibri"
$zh.Range("I3").Font.Size = 11
$zh.Range("I3").Font.Color = 15570276

# ------------------------------------------------------------------
# de-de sheet
# ------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Columns.Item(3).ColumnWidth = 29.17
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17

# "Status" column shares the same text as the Overview status cells.
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 (ac55d992-...)
$de.Range("I2").Value = $fileA
$de.Range("J2").Value = "ac55d992-1edd-410e-bdfd-f7bebed1963e.345e5284197f96f7f8137abb0f1f643e2aedff4e.de-de.xlf"
$de.Range("K2").Value = "2016-08-24 07:03:17"

# Row 3 (c9c70c95-...)
$de.Range("I3").Value = $fileB
$de.Range("J3").Value = "c9c70c95-c479-43fd-a4d2-48270d45b7c0.5ce180b0445c7c4ad289ac1955311cb07e44f947.de-de.xlf"
$de.Range("K3").Value = "2016-08-24 07:03:17"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "$repoBase/$fileA", "", "", $fileA)
$de.Hyperlinks.Add($de.Range("I2"), "$repoBase/$fileA", "", "", $fileA)
$de.Hyperlinks.Add($de.Range("A3"), "$repoBase/$fileB", "", "", $fileB)
$de.Hyperlinks.Add($de.Range("I3"), "$repoBase/$fileB", "", "", $fileB)

$de.Range("I2").Font.Underline = 1
$de.Range("I2").Font.Name = "Calibri"
$de.Range("I2").Font.Size = 11
$de.Range("I2").Font.Color = 15570276
$de.Range("I3").Font.Underline = 1
$de.Range("I3").Font.Name = "Calibri"
$de.Range("I3").Font.Size = 11
$de.Range("I3").Font.Color = 15570276

Write-Host "Handback report generated."
